$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 7951
$ws1.Range("F5").Value = 96
$ws1.Range("F6").Value = 221
$ws1.Range("F9").Value = 113
$ws1.Range("F10").Value = 464
$ws1.Range("F15").Value = 73
$ws1.Range("F17").Value = 5834
$ws1.Range("F18").Value = 178
$ws1.Range("F19").Value = 259
$ws1.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"
$ws1.Range("F20").Value = 1745
$ws1.Range("F21").Value = 236
$ws1.Range("F22").Value = 386

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 7951
$ws4.Range("F5").Value = 96
$ws4.Range("F6").Value = 221
$ws4.Range("F9").Value = 113
$ws4.Range("F10").Value = 464
$ws4.Range("F15").Value = 73
$ws4.Range("F18").Value = 5834
$ws4.Range("F20").Value = 178
$ws4.Range("F21").Value = 259
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"
$ws4.Range("F22").Value = 1745
$ws4.Range("F23").Value = 236
$ws4.Range("F24").Value = 386
